# Updates cryptos list values scraped on Sat Nov 30 03:37:22 UTC 2024 (GitHub Actions run).
# Column D (Price) cells that look like plain decimals are written with a
# text NumberFormat first so Excel keeps them as literal text (e.g. "1.00")
# instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "96.551.56"
# Row 3
$ws.Range("D3").Value = "3.635.95"
$ws.Range("E3").Value = "  +1.38%  "
# Row 4
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.53"
$ws.Range("E5").Value = "  -0.10%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.78"
$ws.Range("E6").Value = "  +14.33%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "653.72"
$ws.Range("E7").Value = "  -0.80%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("E8").Value = "  +3.03%  "
# Row 9
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "
# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.06"
$ws.Range("E10").Value = "  +0.27%  "
# Row 11
$ws.Range("D11").Value = "3.635.84"
$ws.Range("E11").Value = "  +1.41%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.20"
$ws.Range("E12").Value = "  +1.51%  "
# Row 13
$ws.Range("E13").Value = "  +0.23%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.52"
$ws.Range("E14").Value = "  +1.73%  "
# Row 15
$ws.Range("D15").Value = "4.310.62"
$ws.Range("E15").Value = "  +1.28%  "
# Row 16
$ws.Range("D16").Value = "96.395.81"
$ws.Range("E16").Value = "  -0.08%  "
# Row 17
$ws.Range("E17").Value = "  -0.45%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.62"
$ws.Range("E18").Value = "  +11.05%  "
# Row 19
$ws.Range("D19").Value = "3.630.77"
$ws.Range("E19").Value = "  +1.39%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.18"
$ws.Range("E20").Value = "  +4.31%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.37"
$ws.Range("E21").Value = "  +2.67%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.524"
$ws.Range("E22").Value = "  +5.73%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "511.04"
$ws.Range("E23").Value = "  -0.45%  "
# Row 24
$ws.Range("E24").Value = "  -0.46%  "
# Row 25
$ws.Range("E25").Value = "  +1.53%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.85"
$ws.Range("E26").Value = "  -0.12%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "100.80"
$ws.Range("E27").Value = "  +3.95%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.18"
$ws.Range("E28").Value = "  +2.70%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.165"
$ws.Range("E29").Value = "  +10.49%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.03"
$ws.Range("E30").Value = "  -0.03%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.89"
$ws.Range("E31").Value = "  +2.80%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.00%  "
# Row 33
$ws.Range("E33").Value = "  +0.23%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.37"
$ws.Range("E34").Value = "  +5.06%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.34%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.72"
$ws.Range("E36").Value = "  +6.31%  "
# Row 37
$ws.Range("E37").Value = "  +1.66%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.80"
$ws.Range("E38").Value = "  +2.95%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "616.85"
$ws.Range("E39").Value = "  +3.36%  "
# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.63"
$ws.Range("E40").Value = "  +23.17%  "
# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.155"
$ws.Range("E41").Value = "  +2.85%  "
# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.956"
$ws.Range("E42").Value = "  +5.06%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.93"
$ws.Range("E43").Value = "  +4.59%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.14"
$ws.Range("E45").Value = "  +6.23%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0439"
$ws.Range("E46").Value = "  +3.47%  "
# Row 47
$ws.Range("E47").Value = "  +0.18%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.59"
$ws.Range("E48").Value = "  +0.20%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.404"
$ws.Range("E49").Value = "  +13.75%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.53"
$ws.Range("E50").Value = "  +2.88%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.62"
$ws.Range("E51").Value = "  +1.60%  "
